$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Febrero2024"
$ws.Range("B2").Value = "5 Febrero 2024"

$ws.Range("A3").Value = "Marzo2024"
$ws.Range("B3").Value = "9 Marzo 2024"

$ws.Range("A4").Value = "Abril2024"
$ws.Range("B4").Value = "5 Abril 2024"

$ws.Range("A5").Value = "Mayo2024"
$ws.Range("B5").Value = "7 Mayo 2024"
